$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.450.31'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.22%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.179.57'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.40%  '

$ws.Range("E4").Value = '  +0.08%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '602.28'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.41%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '154.70'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.26%  '

$ws.Range("E7").Value = '  +0.05%  '

$ws.Range("B8").Value = 'XRP'
$ws.Range("C8").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.550'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.46%  '

$ws.Range("B9").Value = 'LidoStakedEther'
$ws.Range("C9").Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '3.177.90'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.39%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.158'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.46%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.54'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -10.34%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.512'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.17%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000267'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.05%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '38.64'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.11%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.703.86'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.34%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '66.505.02'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.22%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '7.41'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.35%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.181.06'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.33%  '

$ws.Range("E19").Value = '  +0.45%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '512.67'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.37%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '15.44'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.56%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.732'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.01%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.14'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.45%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '14.85'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.90%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '84.68'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.12%  '

$ws.Range("E26").Value = '  -0.11%  '

$ws.Range("E27").Value = '  -1.52%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.18'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.47%  '

$ws.Range("E29").Value = '  +6.34%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.11'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +7.20%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.98'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.66%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '28.04'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.04%  '

$ws.Range("E33").Value = '  +0.08%  '

$ws.Range("E34").Value = '  -1.86%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.53'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.37%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '514.46'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +5.15%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '54.80'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.55%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0887'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.89%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0420'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.78%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.128'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +6.40%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.83'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.72%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0₃0679'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +4.17%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.299'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.57%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.79'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -8.62%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.44'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.70%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.844.28'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -4.75%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '28.09'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.07%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.38'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.66%  '

$ws.Range("E50").Value = '  +0.34%  '

$ws.Range("E51").Value = '  +5.74%  '
